$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in new test row data (row 5)
$ws.Range("A5").Value = "Testear la lectura de ambos puertos seriales simultaneamente"
$ws.Range("B5").Value = "OK"

# Row 5 height changes to 30 (wraps text, similar to rows 2-3)
$ws.Rows.Item(5).RowHeight = 30

# Update active selection to I6
$ws.Range("I6").Select()
